$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-27 (22 rows, columns A-H) need to be cyclically shifted up by one:
# new row 6  = old row 7
# new row 7  = old row 8
# ...
# new row 26 = old row 27
# new row 27 = old row 6 (with a typo fix in the Discord name: stinkywreslter -> stinkywrestler)

$firstRow = 6
$lastRow = 27

# Capture the original values of row 6 (the row that will wrap around to the bottom)
$savedValues = @()
for ($col = 1; $col -le 8; $col++) {
    $savedValues += $ws.Cells.Item($firstRow, $col).Value2
}

# Shift rows 7..27 up into rows 6..26
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($r - 1, $col).Value = $ws.Cells.Item($r, $col).Value2
    }
}

# Place the saved original row 6 values into row 27, fixing the Discord name typo
for ($col = 1; $col -le 8; $col++) {
    $val = $savedValues[$col - 1]
    if ($col -eq 4) {
        $val = "stinkywrestler#7847"
    }
    $ws.Cells.Item($lastRow, $col).Value = $val
}
